# Generate Report for Handback
# Row 7 ("20168373-6657-4c23-92b2-b8ec96a5e6a6") on both the zh-cn and de-de
# sheets now has a handoff/handback pair recorded, just like rows 2-5 already
# have: a "Latest Target File" hyperlink (column I), a "Latest Handback File"
# name (column J), a "Latest Handback DateTime" (column K) and an
# "Error Detail" message (column P) describing that the handback file isn't
# the most recent one.

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/1505aa6c0f9f7cbe356361adad6495d24aa38ba8/e2e/20168373-6657-4c23-92b2-b8ec96a5e6a6.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c02b23cc71184b582e4981b70225c4e693039dad/e2e/20168373-6657-4c23-92b2-b8ec96a5e6a6.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/1505aa6c0f9f7cbe356361adad6495d24aa38ba8/e2e/20168373-6657-4c23-92b2-b8ec96a5e6a6.md."

# --- zh-cn sheet, row 7 -----------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$iCellZh = $wsZh.Range("I7")
$iCellZh.Value = "20168373-6657-4c23-92b2-b8ec96a5e6a6.md"
$wsZh.Hyperlinks.Add($iCellZh, $latestMdUrl, "", "", "20168373-6657-4c23-92b2-b8ec96a5e6a6.md")

$wsZh.Range("J7").Value = "20168373-6657-4c23-92b2-b8ec96a5e6a6.bbb8e4059c6c68bbf233639570708a7136a61c3d.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-14 03:16:28"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 -------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$iCellDe = $wsDe.Range("I7")
$iCellDe.Value = "20168373-6657-4c23-92b2-b8ec96a5e6a6.md"
$wsDe.Hyperlinks.Add($iCellDe, $latestMdUrl, "", "", "20168373-6657-4c23-92b2-b8ec96a5e6a6.md")

$wsDe.Range("J7").Value = "20168373-6657-4c23-92b2-b8ec96a5e6a6.bbb8e4059c6c68bbf233639570708a7136a61c3d.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-14 03:16:38"
$wsDe.Range("P7").Value = $errorDetail
